# Auto-generated edit script: updates profit/price calculation cells
# across multiple worksheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR)
# per the scheduled runner data refresh.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 973.875
$ws.Range("J129").Value = 1111.8387
$ws.Range("L129").Value = 3335.5161
$ws.Range("N129").Value = -13335.5161
$ws.Range("H137").Value = 1765.2307
$ws.Range("I137").Value = 1464
$ws.Range("J137").Value = 2334.2222
$ws.Range("K137").Value = 4392
$ws.Range("L137").Value = 7002.6666
$ws.Range("M137").Value = -1842
$ws.Range("N137").Value = -12102.6666
$ws.Range("H138").Value = 2503778.2
$ws.Range("I138").Value = 8335177.5
$ws.Range("J138").Value = 4607.232
$ws.Range("K138").Value = 25005532.5
$ws.Range("L138").Value = 13821.696
$ws.Range("M138").Value = -25000392.5
$ws.Range("N138").Value = -24101.696
$ws.Range("H140").Value = 77693.336
$ws.Range("J140").Value = 82717.5
$ws.Range("L140").Value = 82717.5
$ws.Range("N140").Value = -93077.5
$ws.Range("H141").Value = 14124.429
$ws.Range("I141").Value = 7728.5
$ws.Range("K141").Value = 23185.5
$ws.Range("M141").Value = -18005.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1272.5555
$ws.Range("I74").Value = 1121.1765
$ws.Range("K74").Value = 1121.1765
$ws.Range("M74").Value = -247.1765
$ws.Range("H77").Value = 1272.5555
$ws.Range("I77").Value = 1121.1765
$ws.Range("K77").Value = 5605.8825
$ws.Range("M77").Value = -1237.8825
$ws.Range("H132").Value = 1468.2
$ws.Range("I132").Value = 1074.1428
$ws.Range("J132").Value = 2387.6667
$ws.Range("K132").Value = 3222.4284
$ws.Range("L132").Value = 7163.000100000001
$ws.Range("M132").Value = -692.4284000000002
$ws.Range("N132").Value = -12223.0001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 1970.1428
$ws.Range("I134").Value = 1735.2778
$ws.Range("J134").Value = 2620.5386
$ws.Range("K134").Value = 5205.8334
$ws.Range("L134").Value = 7861.6158
$ws.Range("M134").Value = -2670.8334
$ws.Range("N134").Value = -12931.6158

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1639
$ws.Range("I16").Value = 1582.2
$ws.Range("J16").Value = 1733.6666
$ws.Range("K16").Value = 1582.2
$ws.Range("L16").Value = 1733.6666
$ws.Range("M16").Value = -1295.2
$ws.Range("N16").Value = -2307.6666
$ws.Range("H31").Value = 22224936
$ws.Range("I31").Value = 30305052
$ws.Range("J31").Value = 4614.25
$ws.Range("K31").Value = 30305052
$ws.Range("L31").Value = 4614.25
$ws.Range("M31").Value = -30304757
$ws.Range("N31").Value = -5204.25
$ws.Range("H34").Value = 22224936
$ws.Range("I34").Value = 30305052
$ws.Range("J34").Value = 4614.25
$ws.Range("K34").Value = 30305052
$ws.Range("L34").Value = 4614.25
$ws.Range("M34").Value = -30304850
$ws.Range("N34").Value = -5018.25
$ws.Range("H47").Value = 15633.333
$ws.Range("I47").Value = 10000
$ws.Range("J47").Value = 18450
$ws.Range("K47").Value = 10000
$ws.Range("L47").Value = 18450
$ws.Range("M47").Value = -9434
$ws.Range("N47").Value = -19582
$ws.Range("H64").Value = 30000
$ws.Range("J64").Value = 30000
$ws.Range("L64").Value = 30000
$ws.Range("N64").Value = -30496
$ws.Range("H67").Value = 30000
$ws.Range("J67").Value = 30000
$ws.Range("L67").Value = 30000
$ws.Range("N67").Value = -31716
$ws.Range("H113").Value = 1639
$ws.Range("I113").Value = 1582.2
$ws.Range("J113").Value = 1733.6666
$ws.Range("K113").Value = 1582.2
$ws.Range("L113").Value = 1733.6666
$ws.Range("M113").Value = 587.8
$ws.Range("N113").Value = -6073.6666
$ws.Range("H138").Value = 44458.9
$ws.Range("J138").Value = 44458.9
$ws.Range("L138").Value = 44458.9
$ws.Range("N138").Value = -54738.9

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 998.75
$ws.Range("I92").Value = 998.75
$ws.Range("K92").Value = 2996.25
$ws.Range("M92").Value = -1748.25
$ws.Range("H121").Value = 42655.582
$ws.Range("J121").Value = 84638.836
$ws.Range("L121").Value = 253916.508
$ws.Range("N121").Value = -256536.508
$ws.Range("H123").Value = 4066.6667
$ws.Range("I123").Value = 1100
$ws.Range("J123").Value = 10000
$ws.Range("K123").Value = 3300
$ws.Range("L123").Value = 30000
$ws.Range("M123").Value = -850
$ws.Range("N123").Value = -34900
$ws.Range("H131").Value = 865.66
$ws.Range("I131").Value = 487.8
$ws.Range("J131").Value = 885.54736
$ws.Range("K131").Value = 1463.4
$ws.Range("L131").Value = 2656.64208
$ws.Range("M131").Value = 3576.6
$ws.Range("N131").Value = -12736.64208
$ws.Range("H139").Value = 2437.0645
$ws.Range("I139").Value = 1915.8334
$ws.Range("K139").Value = 5747.5002
$ws.Range("M139").Value = -607.5002000000004

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H93").Value = 20250
$ws.Range("J93").Value = 20250
$ws.Range("L93").Value = 20250
$ws.Range("N93").Value = -23994
$ws.Range("H102").Value = 1590.12
$ws.Range("I102").Value = 1531.738
$ws.Range("J102").Value = 1896.625
$ws.Range("K102").Value = 1531.738
$ws.Range("L102").Value = 1896.625
$ws.Range("M102").Value = 90.26199999999994
$ws.Range("N102").Value = -5140.625
$ws.Range("H126").Value = 2179.7812
$ws.Range("I126").Value = 1954.4584
$ws.Range("J126").Value = 2855.75
$ws.Range("K126").Value = 5863.3752
$ws.Range("L126").Value = 8567.25
$ws.Range("M126").Value = -3393.3752
$ws.Range("N126").Value = -13507.25
$ws.Range("H132").Value = 1639.0962
$ws.Range("I132").Value = 1347
$ws.Range("J132").Value = 3034.6667
$ws.Range("K132").Value = 4041
$ws.Range("L132").Value = 9104.000100000001
$ws.Range("M132").Value = -1511
$ws.Range("N132").Value = -14164.0001
$ws.Range("H139").Value = 138134.5
$ws.Range("J139").Value = 138134.5
$ws.Range("L139").Value = 138134.5
$ws.Range("N139").Value = -148414.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3246.5557
$ws.Range("I40").Value = 4263
$ws.Range("J40").Value = 1976
$ws.Range("K40").Value = 4263
$ws.Range("L40").Value = 1976
$ws.Range("M40").Value = -4127
$ws.Range("N40").Value = -2248
$ws.Range("H64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("N64").ClearContents()
$ws.Range("H67").Value = 0
$ws.Range("J67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("N67").ClearContents()
$ws.Range("H122").Value = 8933155
$ws.Range("I122").Value = 13163027
$ws.Range("J122").Value = 3426.6667
$ws.Range("K122").Value = 39489081
$ws.Range("L122").Value = 10280.0001
$ws.Range("M122").Value = -39486631
$ws.Range("N122").Value = -15180.0001
$ws.Range("H133").Value = 56185.547
$ws.Range("J133").Value = 56185.547
$ws.Range("L133").Value = 56185.547
$ws.Range("N133").Value = -61245.547
$ws.Range("H134").Value = 75514.5
$ws.Range("J134").Value = 75514.5
$ws.Range("L134").Value = 75514.5
$ws.Range("N134").Value = -85654.5
$ws.Range("H137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("N137").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H63").Value = 131600
$ws.Range("J63").Value = 131600
$ws.Range("L63").Value = 131600
$ws.Range("N63").Value = -132848
$ws.Range("H66").Value = 131600
$ws.Range("J66").Value = 131600
$ws.Range("L66").Value = 394800
$ws.Range("N66").Value = -401040
$ws.Range("H113").Value = 720.5278
$ws.Range("I113").Value = 513.24
$ws.Range("K113").Value = 1539.72
$ws.Range("M113").Value = 630.28
$ws.Range("H136").Value = 1338.0526
$ws.Range("I136").Value = 1328.2609
$ws.Range("J136").Value = 1353.0667
$ws.Range("K136").Value = 3984.7827
$ws.Range("L136").Value = 4059.2001
$ws.Range("M136").Value = -1434.7827
$ws.Range("N136").Value = -9159.2001
$ws.Range("H138").Value = 72423.336
$ws.Range("J138").Value = 72423.336
$ws.Range("L138").Value = 72423.336
$ws.Range("N138").Value = -82703.336
